# Apply the changes described by the diff:
# 1. In the "day" sheet, cell D19 (bsecode for Colgate Palmolive) is converted
#    from a text value to a true numeric value (500830).
# 2. In the "week" sheet, 11 new data rows (25-35) are appended, extending the
#    sheet from row 24 to row 35. The bsecode column (D) for these new rows is
#    stored as text, matching the existing "inlineStr" pattern used elsewhere
#    in the sheet for newly scraped rows.

$wb = $excel.ActiveWorkbook

# --- 1. Fix D19 on the "day" sheet: text "500830" -> numeric 500830 ---
$dayWs = $wb.Worksheets.Item("day")
$dayWs.Range("D19").Value = 500830

# --- 2. Append new rows 25-35 on the "week" sheet ---
$weekWs = $wb.Worksheets.Item("week")

# New row data: sr, nsecode, name, bsecode, per_chg, close, volume, timeframe, datetime
$newRows = @(
    @(1,  "BAJAJ-AUTO", "Bajaj Auto Limited",                              "532977", 0.39,  9961.75, 320933,    "week", "16/06/2024 11:34:20"),
    @(2,  "COFORGE",    "Coforge (Niit Tech)",                             "532541", -1.36, 5201.75, 354469,    "week", "16/06/2024 11:34:20"),
    @(3,  "PIDILITIND", "Pidilite Industries Limited",                     "500331", 1.01,  3109.8,  366384,    "week", "16/06/2024 11:34:20"),
    @(4,  "ASIANPAINT", "Asian Paints Limited",                            "500820", 0.4,   2921.6,  982815,    "week", "16/06/2024 11:34:20"),
    @(5,  "DEEPAKNTR",  "Deepak Nitrite Limited",                          "506401", 0.11,  2413.7,  450097,    "week", "16/06/2024 11:34:20"),
    @(6,  "HDFCBANK",   "Hdfc Bank Limited",                               "500180", 1.02,  1596.9,  12770277,  "week", "16/06/2024 11:34:20"),
    @(7,  "CIPLA",      "Cipla Limited",                                   "500087", 1.31,  1564.75, 1845043,   "week", "16/06/2024 11:34:20"),
    @(8,  "TATAMOTORS", "Tata Motors Limited",                             "500570", 0.77,  993.4,   11591421,  "week", "16/06/2024 11:34:20"),
    @(9,  "PEL",        "Piramal Enterprises Limited",                     "500302", 0.22,  882.5,   1247110,   "week", "16/06/2024 11:34:20"),
    @(10, "ITC",        "Itc Limited",                                     "500875", 0.2,   431.15,  9217804,   "week", "16/06/2024 11:34:20"),
    @(11, "M&MFIN",     "Mahindra & Mahindra Financial Services Limited",  "532720", 1.93,  298.95,  4666536,   "week", "16/06/2024 11:34:20")
)

$startRow = 25
$endRow = $startRow + $newRows.Length - 1

# Pre-format column D for the new block as Text so the numeric-looking
# bsecode strings are stored as text (t="inlineStr"/"s"), not numbers -
# matching the rest of the bsecode column for freshly appended batches.
$colDRangeAddress = "D" + $startRow + ":D" + $endRow
$weekWs.Range($colDRangeAddress).NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $weekWs.Cells.Item($r, 1).Value = $row[0]   # A: sr
    $weekWs.Cells.Item($r, 2).Value = $row[1]   # B: nsecode
    $weekWs.Cells.Item($r, 3).Value = $row[2]   # C: name
    $weekWs.Cells.Item($r, 4).Value = $row[3]   # D: bsecode (text)
    $weekWs.Cells.Item($r, 5).Value = $row[4]   # E: per_chg
    $weekWs.Cells.Item($r, 6).Value = $row[5]   # F: close
    $weekWs.Cells.Item($r, 7).Value = $row[6]   # G: volume
    $weekWs.Cells.Item($r, 8).Value = $row[7]   # H: timeframe
    $weekWs.Cells.Item($r, 9).Value = $row[8]   # I: Date Time
}

Write-Host "Done."
